$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B17").Value = "Mahalakshmi"
$ws.Range("B18").Value = "Mahalakshmi"
$ws.Range("B19").Value = "Mahalakshmi"
$ws.Range("B20").Value = "Mahalakshmi"

$ws.Range("B16:B20").Select()
